$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.286999999999999
$ws.Range("B21").Value = 6.459000000000001
$ws.Range("B23").Value = 6.842000000000001
$ws.Range("B25").Value = 6.556999999999999
